# Regenerate the "K" column (column G) values for the save_data sheet.
# The commit replaces the old "Strike#" derived values with the real
# strikeout count "K" for each outing (row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), per the regenerated save data.
$kValues = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 3
    11 = 1
    12 = 1
    13 = 1
    15 = 0
    16 = 2
    17 = 3
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 2
    37 = 3
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    48 = 1
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    56 = 1
    57 = 2
    58 = 1
    59 = 0
    60 = 1
    61 = 0
    62 = 0
    63 = 1
    64 = 1
    65 = 1
    67 = 2
    68 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
